{"js": "// Update the DaCapo ZGC graphchi (heap-2G) benchmark results table.\n// The document is a single-column table where each row holds one stat\n// value. A handful of scalar values changed, and the three \"raw timing\n// line\" rows (originally multi-run tab-separated text) were collapsed\n// down to the single summary number that the corresponding header rows\n// (0, 1, 2) used to hold.\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"items\");\nawait context.sync();\n\n// Map of row index (0-based) -> new cell text.\nconst updates = {\n  0: \"0M\",\n  1: \"0M\",\n  2: \"0M\",\n  3: \"170\",\n  6: \"0.01104\",\n  7: \"0.00206\",\n  11: \"0.71122\",\n  43: \"98.53\",\n  44: \"0.71\",\n  45: \"48\",\n};\n\nfor (const [rowIndex, newText] of Object.entries(updates)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Update the DaCapo ZGC graphchi (heap-2G) benchmark results table.\n# The document is a single-column table where each row holds one stat\n# value. A handful of scalar values changed, and the three \"raw timing\n# line\" rows (originally multi-run tab-separated text) were collapsed\n# down to the single summary number that the corresponding header rows\n# (1, 2, 3 in 1-based COM indexing) used to hold.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Map of 1-based row number -> new cell text (COM Cell/Rows are 1-indexed).\n$updates = [ordered]@{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"170\"\n    7  = \"0.01104\"\n    8  = \"0.00206\"\n    12 = \"0.71122\"\n    44 = \"98.53\"\n    45 = \"0.71\"\n    46 = \"48\"\n}\n\nforeach ($rowNum in $updates.Keys) {\n    $cell = $t.Cell($rowNum, 1)\n    $cell.Range.Text = $updates[$rowNum]\n}\n"}
